$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $found = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $find"
    } else {
        # Assign directly to Range.Text (instead of using Find's Replace
        # parameter) so that straight apostrophes are not auto-converted
        # into curly "smart quote" apostrophes by the replace engine.
        $range.Text = $replace
    }
}

Replace-Text "Empty group means one." "Un groupe vide signifie un."
Replace-Text "Also, try starting with " "Essayez également de commencer par "
Replace-Text "and then decrease the last number gradually." "puis diminuez progressivement le dernier chiffre."
Replace-Text "The table above represents what we know thanks to the answers given by the second mathematician (Fil)." "Le tableau ci-dessus représente ce que l'on sait grâce aux réponses données par le deuxième mathématicien (Fil)."
Replace-Text "But another important information is that the first mathematician (Mike) is not able to know the correct combination, " "Mais une autre information importante est que le premier mathématicien (Mike) n'est pas capable de connaître la bonne combinaison, "
Replace-Text "even if he knows the actual value of the sum!" "même s'il connaît la valeur réelle de la somme !"
Replace-Text "This can only happen if the number corresponding to the correct sum appears more than once in the list! (otherwise he would have guessed the correct numbers after the second question) So, the sum must be 13, and the corresponding combinations are:" "Cela ne peut arriver que si le nombre correspondant à la bonne somme apparaît plus d’une fois dans la liste ! (sinon il aurait deviné les bons nombres après la deuxième question) Donc, la somme doit être 13, et les combinaisons correspondantes sont :"
Replace-Text "The final clue is that the youngest child has blue eyes." "Le dernier indice est que le plus jeune enfant a les yeux bleus."
Replace-Text "What we get from this clue is that now we know that a youngest child exists!" "Ce que nous retenons de cet indice, c'est que nous savons désormais qu'un plus jeune enfant existe !"
Replace-Text "So " "Donc "
Replace-Text "is not possible and " "n'est pas possible et "
Replace-Text "is the only remaining option." "est la seule option restante."
